# ProjectProposal.docx -> "Roguelike Project Outline" edit
#
# Summary of changes (see commit message / diff):
#  - Title: "Roguelike Project Proposal" -> "Roguelike Project Outline";
#    drop the spell-check proofErr wrapper around "Roguelike"; the
#    _GoBack bookmark moves from the end of the document to the title.
#  - Team member lines: collapse the "First" + "Last" run pairs (that
#    were wrapped in spellStart/spellEnd proofErr markers) into single
#    runs per line.
#  - Sprint 4 paragraph: drop the gramStart/gramEnd proofErr markers and
#    merge the surrounding runs.
#  - Final paragraph: add underline to the paragraph mark's run
#    properties and drop the trailing _GoBack bookmark (it moved to the
#    title).
#  - Append a page break plus a new "Class Design" section describing
#    the Game / Input / Tile / Floor / GameObject / Character / Player /
#    Enemy / Item / Inventory classes.

# NOTE: this interpreter loses COM object identity when a COM object is
# bound to a *named* function parameter (e.g. `-Paragraph $p`); plain
# positional parameter binding works fine, so every helper below takes
# its arguments positionally.

$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($Paragraph, $InnerXml) {
    $range = $Paragraph.Range
    $xml = "<w:p $W>$InnerXml</w:p>"
    [void]$range.InsertXML($xml)
}

function Insert-ParagraphsAfter($Paragraph, $PsXml) {
    # NB: Range.Collapse(...) is a no-op in this host, so build a fresh
    # zero-length range sitting exactly at the paragraph's end instead
    # of trying to collapse the paragraph's own Range in place.
    $range = $Paragraph.Range
    $insertPoint = $d.Range($range.End, $range.End)
    [void]$insertPoint.InsertXML($PsXml)
}

# ---------------------------------------------------------------------
# 1) Title paragraph: merge runs, drop proofErr, retitle, move bookmark
# ---------------------------------------------------------------------
$titlePPr = '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr>'
$titleRunProps = '<w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr>'
$titleInner = $titlePPr +
    "<w:r>$titleRunProps<w:t>R</w:t></w:r>" +
    "<w:r>$titleRunProps<w:t xml:space=`"preserve`">oguelike Project </w:t></w:r>" +
    "<w:r>$titleRunProps<w:t>Outline</w:t></w:r>" +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-ParagraphXml $d.Paragraphs.Item(1) $titleInner

# ---------------------------------------------------------------------
# 2) Team members block (paragraph 3): merge "Ioan" + " " + "Istrate"
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$memberProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$memberUProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr>'
$p3Inner = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    "<w:r>$memberUProps<w:t>Project Manager</w:t></w:r>" +
    "<w:r>$memberProps<w:t>:</w:t></w:r>" +
    "<w:r>$memberProps<w:tab/></w:r>" +
    "<w:r>$memberProps<w:t>Derek Royse</w:t></w:r>" +
    "<w:r>$memberProps<w:br/></w:r>" +
    "<w:r>$memberUProps<w:t>Team Members</w:t></w:r>" +
    "<w:r>$memberProps<w:t>:</w:t></w:r>" +
    "<w:r>$memberProps<w:tab/></w:r>" +
    "<w:r>$memberProps<w:t>Ioan Istrate</w:t></w:r>"
Set-ParagraphXml $p3 $p3Inner

# ---------------------------------------------------------------------
# 3) Gary Danovich / Adam McCroskey / Andy Pritt: merge first+last runs
# ---------------------------------------------------------------------
function Set-NameParagraph($Paragraph, $FullName) {
    $indPPr = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="1440" w:firstLine="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'
    $inner = $indPPr + "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>$FullName</w:t></w:r>"
    Set-ParagraphXml $Paragraph $inner
}

Set-NameParagraph $d.Paragraphs.Item(5) "Gary Danovich"
Set-NameParagraph $d.Paragraphs.Item(6) "Adam McCroskey"
Set-NameParagraph $d.Paragraphs.Item(8) "Andy Pritt"

# ---------------------------------------------------------------------
# 4) Sprint 4 paragraph (21): drop gramStart/gramEnd, merge runs
# ---------------------------------------------------------------------
$p21 = $d.Paragraphs.Item(21)
$sProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$sUProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr>'
$sUSupProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:vertAlign w:val="superscript"/></w:rPr>'
$p21Inner = '<w:pPr><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    "<w:r>$sProps<w:t>T</w:t></w:r>" +
    "<w:r>$sProps<w:t xml:space=`"preserve`">he team will expand the combat, enemy, inventory, and leveling systems in Sprint </w:t></w:r>" +
    "<w:r>$sProps<w:t>4</w:t></w:r>" +
    "<w:r>$sProps<w:t>. This</w:t></w:r>" +
    "<w:r>$sProps<w:t xml:space=`"preserve`"> sprint will also the see the completion and refinement of the input and output systems, and will be finished on </w:t></w:r>" +
    "<w:r>$sUProps<w:t>May 9</w:t></w:r>" +
    "<w:r>$sUSupProps<w:t>th</w:t></w:r>" +
    "<w:r>$sUProps<w:t>.</w:t></w:r>"
Set-ParagraphXml $p21 $p21Inner

# ---------------------------------------------------------------------
# 5) Final paragraph (22): add underline to pPr/rPr, drop the bookmark
#    (it was moved up to the title paragraph in step 1)
# ---------------------------------------------------------------------
$p22 = $d.Paragraphs.Item(22)
$fProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$fUProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr>'
$fUSupProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:vertAlign w:val="superscript"/></w:rPr>'
$p22Inner = '<w:pPr><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' +
    "<w:r>$fProps<w:t>In the final sprint, the team will focus</w:t></w:r>" +
    "<w:r>$fProps<w:t xml:space=`"preserve`"> on advanced combat and a player class system. These enhancements would be in place for the projects final version to be presented on </w:t></w:r>" +
    "<w:r>$fUProps<w:t>May 16</w:t></w:r>" +
    "<w:r>$fUSupProps<w:t>th</w:t></w:r>" +
    "<w:r>$fUProps<w:t>.</w:t></w:r>"
Set-ParagraphXml $p22 $p22Inner

# ---------------------------------------------------------------------
# 6) Append: page break paragraph, "Class Design" header, and the new
#    class description paragraphs
# ---------------------------------------------------------------------
$u24 = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr>'
$pageBreakPara = "<w:p $W><w:pPr>$u24</w:pPr><w:r>$u24<w:br w:type=`"page`"/></w:r></w:p>"

$headerProps = '<w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr>'
$headerPara = "<w:p $W><w:pPr>$headerProps</w:pPr><w:r>$headerProps<w:lastRenderedPageBreak/><w:t>Class Design</w:t></w:r></w:p>"

$indPPr = '<w:pPr><w:ind w:left="2880" w:hanging="2880"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'
$bProps = '<w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rProps = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$gamePara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Game</w:t></w:r>" +
    "<w:r>$rProps<w:t xml:space=`"preserve`">  </w:t></w:r>" +
    "<w:r>$rProps<w:tab/><w:t>S</w:t></w:r>" +
    "<w:r>$rProps<w:t>tarts/ends the game. May perform other large functions between objects.</w:t></w:r>" +
    "<w:r>$rProps<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r>$rProps<w:t>The heart of the program.</w:t></w:r>" +
    "</w:p>"

$inputPara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Input</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t xml:space=`"preserve`">Receives input from the player and translates it to a game </w:t></w:r>" +
    "<w:r>$rProps<w:t>f</w:t></w:r>" +
    "<w:r>$rProps<w:t>unction.</w:t></w:r>" +
    "</w:p>"

$tilePara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Tile</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>An indi</w:t></w:r>" +
    "<w:r>$rProps<w:t>vidual space in the game world.</w:t></w:r>" +
    "<w:r>$rProps<w:br/></w:r>" +
    "<w:r>$rProps<w:t>Tiles can contain ground, walls, empty spaces, the player, enemies, items, stairs, etc.</w:t></w:r>" +
    "</w:p>"

$floorPara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Floor</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>Floor is an object made up of Tiles. Represents one level of the dungeon.</w:t></w:r>" +
    "</w:p>"

$gameObjectPara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>GameObject</w:t></w:r>" +
    "<w:r>$bProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>Simple base class that will contain coordinates and very simple methods</w:t></w:r>" +
    "<w:r>$rProps<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r>$rProps<w:t>that apply to all game objects.</w:t></w:r>" +
    "</w:p>"

$characterPara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Character</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>Base class for player and enemies. Will ideally contain movement and</w:t></w:r>" +
    "<w:r>$rProps<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r>$rProps<w:t>attack methods that can be inherited by its children.</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "</w:p>"

$playerPara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Player</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>Child of Game Object/Character. It's the player, and will include attributes</w:t></w:r>" +
    "<w:r>$rProps<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r>$rProps<w:t>related to combat and leveling.</w:t></w:r>" +
    "</w:p>"

$enemyPara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Enemy</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>Child of Object/Character. It can move around and attack the player.</w:t></w:r>" +
    "</w:p>"

$itemPara = "<w:p $W>$indPPr" +
    "<w:r>$bProps<w:t>Item</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>Base class for objects the player ca</w:t></w:r>" +
    "<w:r>$rProps<w:t xml:space=`"preserve`">n pick up to heal/modify stats. </w:t></w:r>" +
    "</w:p>"

$plainPPr = '<w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'
$inventoryPara = "<w:p $W>$plainPPr" +
    "<w:r>$bProps<w:t>Inventory</w:t></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:tab/></w:r>" +
    "<w:r>$rProps<w:t>An object that Characters have that will hold Item objects.</w:t></w:r>" +
    "</w:p>"

$allNewParas = $pageBreakPara + $headerPara + $gamePara + $inputPara + $tilePara + $floorPara +
    $gameObjectPara + $characterPara + $playerPara + $enemyPara + $itemPara + $inventoryPara

Insert-ParagraphsAfter $d.Paragraphs.Item(22) $allNewParas

Write-Output "edit complete"
